{"js": "// Apply the documented edit: add a new \"Step 10\" paragraph about saving the\n// Easy Localizer excel as .xlsx, and append a new \"Note\" section (with fix\n// steps for the .NET 4 resource-file error) at the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1) Locate the \"Step 9: ...\" paragraph and insert the new \"Step 10\"\n//        paragraph right after it (and before the \"Note: - \" paragraph).\nconst step9 = paragraphs.items.find(p => p.text.indexOf(\"Step 9:\") === 0);\nif (!step9) {\n  throw new Error('Could not find the \"Step 9:\" paragraph');\n}\n\nconst step10Text =\n  \"Step 10:-  As excel created by Easy Localizer is in xls format so use \\u201CSave as\\u201D to save excel in xlsx format.\";\nstep9.insertParagraph(step10Text, Word.InsertLocation.after);\nawait context.sync();\n\n// --- 2) Append the new trailing content after the very last paragraph\n//        (\"Step 3: -  Copy the English terms from html file ...\").\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items/text\");\nawait context.sync();\nconst lastParagraph = allParagraphs.items[allParagraphs.items.length - 1];\n\n// Blank paragraph separator.\nlet cursor = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\nconst noteParagraph = cursor.insertParagraph(\n  \"Note: It is quite possible that when you press next in Step 6 you will get error. This error normally comes because Easy reader do not support .NET 4 resource files.\",\n  Word.InsertLocation.after\n);\n\ncursor = noteParagraph.insertParagraph(\n  \"To fix this error follow following steps:-\",\n  Word.InsertLocation.after\n);\n\ncursor = cursor.insertParagraph(\n  \"Step 1:-  Open the resource file in notepad++ which is causing error in Easy Reader.\",\n  Word.InsertLocation.after\n);\n\ncursor = cursor.insertParagraph(\n  \"Step 2: Use replace all feature of notepad++ to replace all instances of \\\"Version=4.0.0.0\\\" with \\\"Version=2.0.0.0\\\". Do this procedure with all the resource files giving error.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\n// Now that all paragraphs exist, bold only the \"Note:\" lead-in of the note\n// paragraph (doing this before creating the later paragraphs would leak the\n// bold formatting into them).\nconst boldRange = noteParagraph.search(\"Note:\", { matchCase: true });\nboldRange.load(\"items\");\nawait context.sync();\nboldRange.items[0].font.set({ bold: true });\nawait context.sync();\n", "ps1": "# Apply the documented edit: add a new \"Step 10\" paragraph about saving the\n# Easy Localizer excel as .xlsx, and append a new \"Note\" section (with fix\n# steps for the .NET 4 resource-file error) at the end of the document.\n\n$d = $word.ActiveDocument\n\n# --- 1) Locate the \"Step 9: ...\" paragraph and insert the new \"Step 10\"\n#        paragraph right after it (and before the \"Note: - \" paragraph).\n$step9Index = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.IndexOf(\"Step 9:\") -eq 0) {\n        $step9Index = $i\n        break\n    }\n}\nif ($step9Index -eq -1) {\n    throw \"Could not find the 'Step 9:' paragraph\"\n}\n\n$step9 = $d.Paragraphs.Item($step9Index)\n$step9.Range.InsertParagraphAfter()\n$step10 = $d.Paragraphs.Item($step9Index + 1)\n$step10.Range.Text = \"Step 10:-  As excel created by Easy Localizer is in xls format so use \" + [char]8220 + \"Save as\" + [char]8221 + \" to save excel in xlsx format.\"\n\n# --- 2) Append the new trailing content after the very last paragraph\n#        (\"Step 3: -  Copy the English terms from html file ...\").\n$last = $d.Paragraphs.Item($d.Paragraphs.Count)\n$last.Range.InsertParagraphAfter()\n\n# Blank paragraph separator.\n$blank = $d.Paragraphs.Item($d.Paragraphs.Count)\n$blank.Range.InsertParagraphAfter()\n\n# Bold \"Note:\" lead-in followed by plain explanatory text.\n$noteParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$noteParagraph.Range.Text = \"Note: It is quite possible that when you press next in Step 6 you will get error. This error normally comes because Easy reader do not support .NET 4 resource files.\"\n$noteParagraph.Range.InsertParagraphAfter()\n\n$toFixParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$toFixParagraph.Range.Text = \"To fix this error follow following steps:-\"\n$toFixParagraph.Range.InsertParagraphAfter()\n\n$fixStep1Paragraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$fixStep1Paragraph.Range.Text = \"Step 1:-  Open the resource file in notepad++ which is causing error in Easy Reader.\"\n$fixStep1Paragraph.Range.InsertParagraphAfter()\n\n$fixStep2Paragraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$fixStep2Paragraph.Range.Text = 'Step 2: Use replace all feature of notepad++ to replace all instances of \"Version=4.0.0.0\" with \"Version=2.0.0.0\". Do this procedure with all the resource files giving error.'\n\n# Bold only the \"Note:\" lead-in (first 5 characters) of the note paragraph.\n# Re-fetch the paragraph's range start since InsertParagraphAfter calls above\n# may have shifted later ranges.\n$noteRangeStart = $noteParagraph.Range.Start\n$boldRange = $d.Range($noteRangeStart, $noteRangeStart + 5)\n$boldRange.Font.Bold = 1\n"}
